# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 5 (Fecha 45037 = 2023-04-21),
# pushing the existing rows 5-11 down to rows 6-12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(5).Insert()

$ws.Cells.Item(5, 1).Value  = 7
$ws.Cells.Item(5, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(5, 3).Value  = "Ñuble"
$ws.Cells.Item(5, 4).Value  = 45037
$ws.Cells.Item(5, 5).Value  = 16
$ws.Cells.Item(5, 6).Value  = "Fruta"
$ws.Cells.Item(5, 7).Value  = 100107
$ws.Cells.Item(5, 8).Value  = "Otros"
$ws.Cells.Item(5, 9).Value  = 100107011
$ws.Cells.Item(5, 10).Value = "Tuna"
$ws.Cells.Item(5, 11).Value = "Sin especificar"
$ws.Cells.Item(5, 12).Value = "Primera"
$ws.Cells.Item(5, 13).Value = 60
$ws.Cells.Item(5, 14).Value = 16000
$ws.Cells.Item(5, 15).Value = 16000
$ws.Cells.Item(5, 16).Value = 16000
$ws.Cells.Item(5, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(5, 18).Value = "Región Metropolitana"
$ws.Cells.Item(5, 19).Value = 889
$ws.Cells.Item(5, 20).Value = 18
